$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.412.59'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.46%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.867.83'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.92%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.98%  '

$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4684'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.73%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3970'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.52%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.55'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.93%  '

$ws.Range("E10").Value = '  +1.80%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9989'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.75%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.97'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.045'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.68%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.861.89'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.20%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.266'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.90%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.66'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.96%  '

$ws.Range("E17").Value = '  +0.02%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001040'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.87%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06622'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.13%  '

$ws.Range("E21").Value = '  +0.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '28.430.22'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.475'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.48%  '

$ws.Range("E24").Value = '  +2.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.271'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.089.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.24%  '

$ws.Range("E28").Value = '  +1.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.118'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.62%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.485'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.75%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '120.06'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9671'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.19%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09508'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.38%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.589'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.36%  '

$ws.Range("E35").Value = '  +2.17%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.373'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06095'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.57%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02253'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.39%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.334'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.38%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.177'
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5934'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.32%  '

$ws.Range("E42").Value = '  +0.08%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1873'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.69%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.279'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.40%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5569'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.42%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.12'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.80%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.955'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07164'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.062'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +12.81%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '111.76'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.19%  '
